# 23 dec 2023 update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MD10000.1-OCT")

# Fill in the next 4 days of loan payments (rows 4-7, columns J:L)
$ws.Range("J4").Value = 45275
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 1

$ws.Range("J5").Value = 45276
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 1

$ws.Range("J6").Value = 45277
$ws.Range("K6").Value = 100
$ws.Range("L6").Value = 1

$ws.Range("J7").Value = 45278
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 1

# Make sure the date style used for J3 carries over to the new J4:J7 cells
$ws.Range("J3").Copy()
$ws.Range("J4:J7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update selection to reflect where the user left off editing
$ws.Range("L3:L7").Select()
